$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates: force text to preserve exact formatting ---
$ws.Range("D2").Value = "'30.817.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'1.892.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'0.9994"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'249.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'0.4767"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.2942"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.06549"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'22.10"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07770"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'97.47"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'1.890.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.7390"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'5.254"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'284.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'30.844.81"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'13.23"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.000007596"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'1.0000"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'2.137.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'5.345"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'0.9992"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'6.266"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'9.271"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'164.49"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'19.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Value = "'1.934"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").Value = "'0.09749"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Value = "'4.322"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Value = "'4.201"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Value = "'0.04891"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = "'1.131"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Value = "'0.7012"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Value = "'0.01917"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Value = "'6.367"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'76.02"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'2.037"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.4287"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'1.000"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.8399"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'102.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'9.424"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'7.109"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'35.87"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'929.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'0.05767"
$ws.Range("D51").Style = "Normal"

# --- Volume(1h) percentage (column E) updates ---
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("E8").Value = "  +0.71%  "
$ws.Range("E9").Value = "  +0.24%  "
$ws.Range("E10").Value = "  +0.23%  "
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("E14").Value = "  -0.45%  "
$ws.Range("E15").Value = "  +1.86%  "
$ws.Range("E16").Value = "  +3.48%  "
$ws.Range("E17").Value = "  +0.90%  "
$ws.Range("E18").Value = "  -1.91%  "
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("E21").Value = "  +0.38%  "
$ws.Range("E22").Value = "  +1.42%  "
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("E24").Value = "  +1.06%  "
$ws.Range("E25").Value = "  -0.86%  "
$ws.Range("E26").Value = "  +0.66%  "
$ws.Range("E27").Value = "  +0.22%  "
$ws.Range("E28").Value = "  -0.81%  "
$ws.Range("E29").Value = "  -1.85%  "
$ws.Range("E30").Value = "  -2.39%  "
$ws.Range("E31").Value = "  -0.89%  "
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("E33").Value = "  +1.86%  "
$ws.Range("E34").Value = "  +1.74%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("E38").Value = "  +2.43%  "
$ws.Range("E39").Value = "  +2.04%  "
$ws.Range("E40").Value = "  +0.46%  "
$ws.Range("E41").Value = "  +6.05%  "
$ws.Range("E42").Value = "  +2.46%  "
$ws.Range("E43").Value = "  +1.22%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("E46").Value = "  -0.79%  "
$ws.Range("E47").Value = "  +0.92%  "
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("E49").Value = "  +0.39%  "
$ws.Range("E50").Value = "  +1.24%  "
$ws.Range("E51").Value = "  +2.16%  "
